$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new rows of data need to be inserted right above the current row 236,
# pushing all existing rows (236..313) down by two (to 238..315).
$ws.Rows.Item(236).Insert()
$ws.Rows.Item(236).Insert()

# Populate the first inserted row (new row 236).
$ws.Cells.Item(236, 1).Value = 3
$ws.Cells.Item(236, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(236, 3).Value = "Coquimbo"
$ws.Cells.Item(236, 4).Value = 44559
$ws.Cells.Item(236, 5).Value = 5
$ws.Cells.Item(236, 6).Value = 100112003
$ws.Cells.Item(236, 7).Value = "Ajo"
$ws.Cells.Item(236, 8).Value = "Chino"
$ws.Cells.Item(236, 9).Value = "1a (cosecha)"
$ws.Cells.Item(236, 10).Value = 65
$ws.Cells.Item(236, 11).Value = 15000
$ws.Cells.Item(236, 12).Value = 16000
$ws.Cells.Item(236, 13).Value = 15462
$ws.Cells.Item(236, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(236, 15).Value = "Llay Llay"
$ws.Cells.Item(236, 16).Value = 1546
$ws.Cells.Item(236, 17).Value = 10
$ws.Cells.Item(236, 18).Value = "Hortaliza"

# Populate the second inserted row (new row 237).
$ws.Cells.Item(237, 1).Value = 3
$ws.Cells.Item(237, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(237, 3).Value = "Coquimbo"
$ws.Cells.Item(237, 4).Value = 44559
$ws.Cells.Item(237, 5).Value = 5
$ws.Cells.Item(237, 6).Value = 100112003
$ws.Cells.Item(237, 7).Value = "Ajo"
$ws.Cells.Item(237, 8).Value = "Chino"
$ws.Cells.Item(237, 9).Value = "2a (cosecha)"
$ws.Cells.Item(237, 10).Value = 48
$ws.Cells.Item(237, 11).Value = 13000
$ws.Cells.Item(237, 12).Value = 13000
$ws.Cells.Item(237, 13).Value = 13000
$ws.Cells.Item(237, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(237, 15).Value = "Llay Llay"
$ws.Cells.Item(237, 16).Value = 1300
$ws.Cells.Item(237, 17).Value = 10
$ws.Cells.Item(237, 18).Value = "Hortaliza"
